$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EmployeeDataBatch16")

# Update the photograph path for all three data rows (D2:D4) to the new path
$newPhotoPath = "C:\Users\magre\IdeaProjects\HRMSB16\src\test\resources\testdata\Batch16.jpg"
$ws.Range("D2").Value = $newPhotoPath
$ws.Range("D3").Value = $newPhotoPath
$ws.Range("D4").Value = $newPhotoPath

# Remove the custom (black) font style that was applied to D3/D4 so they
# revert back to the default cell style used elsewhere in the sheet.
$ws.Range("D2:D4").ClearFormats()

# Update the username values
$ws.Range("E2").Value = "gray000"
$ws.Range("E3").Value = "sanny111"
$ws.Range("E4").Value = "moon222"

# Widen column D to fit the new longer path text (closest value reachable
# through the pixel-quantized ColumnWidth property to the authored 69.109375).
$ws.Columns.Item(4).ColumnWidth = 68.3333333

# Change the active selection from D4 to E4
$ws.Range("E4").Select()
